$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old table data (B2:E5)
$ws.Range("B2:E5").Clear()

# Set the new single cell value
$ws.Range("A2").Value = "jj knb "

# Update the selection to match the new content
$ws.Range("A2").Select()
